$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 122: Everton de Vina vs Deportes Limache
Set-TextCell 122 1 "2025-07-18"
Set-TextCell 122 2 "Everton de Vina"
Set-TextCell 122 3 "Deportes Limache"
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 0
$ws.Cells.Item(122, 6).Value = 1339129
$ws.Cells.Item(122, 7).Value = 6
$ws.Cells.Item(122, 8).Value = 5
$ws.Cells.Item(122, 9).Value = 1
$ws.Cells.Item(122, 10).Value = 2
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 0
$ws.Cells.Item(122, 14).Value = 0
$ws.Cells.Item(122, 15).Value = 0
$ws.Cells.Item(122, 16).Value = 0
Set-TextCell 122 17 "63%"
Set-TextCell 122 18 "37%"
Set-TextCell 122 19 "E"

# Row 123: Union Espanola vs Union La Calera
Set-TextCell 123 1 "2025-07-19"
Set-TextCell 123 2 "Union Espanola"
Set-TextCell 123 3 "Union La Calera"
$ws.Cells.Item(123, 4).Value = 3
$ws.Cells.Item(123, 5).Value = 1
$ws.Cells.Item(123, 6).Value = 1339131
$ws.Cells.Item(123, 7).Value = 10
$ws.Cells.Item(123, 8).Value = 2
$ws.Cells.Item(123, 9).Value = 1
$ws.Cells.Item(123, 10).Value = 4
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 2
$ws.Cells.Item(123, 13).Value = 0
$ws.Cells.Item(123, 14).Value = 0
$ws.Cells.Item(123, 15).Value = 3
$ws.Cells.Item(123, 16).Value = 1
Set-TextCell 123 17 "46%"
Set-TextCell 123 18 "54%"
Set-TextCell 123 19 "L"

# Row 124: Colo Colo vs D. La Serena
Set-TextCell 124 1 "2025-07-19"
Set-TextCell 124 2 "Colo Colo"
Set-TextCell 124 3 "D. La Serena"
$ws.Cells.Item(124, 4).Value = 2
$ws.Cells.Item(124, 5).Value = 1
$ws.Cells.Item(124, 6).Value = 1339126
$ws.Cells.Item(124, 7).Value = 7
$ws.Cells.Item(124, 8).Value = 4
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = 3
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 13).Value = 0
$ws.Cells.Item(124, 14).Value = 0
$ws.Cells.Item(124, 15).Value = 2
$ws.Cells.Item(124, 16).Value = 1
Set-TextCell 124 17 "61%"
Set-TextCell 124 18 "39%"
Set-TextCell 124 19 "L"

# Row 125: Coquimbo Unido vs Deportes Iquique
Set-TextCell 125 1 "2025-07-19"
Set-TextCell 125 2 "Coquimbo Unido"
Set-TextCell 125 3 "Deportes Iquique"
$ws.Cells.Item(125, 4).Value = 4
$ws.Cells.Item(125, 5).Value = 1
$ws.Cells.Item(125, 6).Value = 1339128
$ws.Cells.Item(125, 7).Value = 6
$ws.Cells.Item(125, 8).Value = 2
$ws.Cells.Item(125, 9).Value = 4
$ws.Cells.Item(125, 10).Value = 3
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = 0
$ws.Cells.Item(125, 14).Value = 0
$ws.Cells.Item(125, 15).Value = 4
$ws.Cells.Item(125, 16).Value = 1
Set-TextCell 125 17 "46%"
Set-TextCell 125 18 "54%"
Set-TextCell 125 19 "L"

# Row 126: A. Italiano vs U. Catolica
Set-TextCell 126 1 "2025-07-20"
Set-TextCell 126 2 "A. Italiano"
Set-TextCell 126 3 "U. Catolica"
$ws.Cells.Item(126, 4).Value = 1
$ws.Cells.Item(126, 5).Value = 1
$ws.Cells.Item(126, 6).Value = 1339130
$ws.Cells.Item(126, 7).Value = 2
$ws.Cells.Item(126, 8).Value = 5
$ws.Cells.Item(126, 9).Value = 3
$ws.Cells.Item(126, 10).Value = 5
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 1
$ws.Cells.Item(126, 13).Value = 0
$ws.Cells.Item(126, 14).Value = 0
$ws.Cells.Item(126, 15).Value = 1
$ws.Cells.Item(126, 16).Value = 1
Set-TextCell 126 17 "50%"
Set-TextCell 126 18 "50%"
Set-TextCell 126 19 "E"

# Row 127: Nublense vs Universidad de Chile
Set-TextCell 127 1 "2025-07-20"
Set-TextCell 127 2 "Nublense"
Set-TextCell 127 3 "Universidad de Chile"
$ws.Cells.Item(127, 4).Value = 2
$ws.Cells.Item(127, 5).Value = 2
$ws.Cells.Item(127, 6).Value = 1339132
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 3
$ws.Cells.Item(127, 9).Value = 1
$ws.Cells.Item(127, 10).Value = 4
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = 0
$ws.Cells.Item(127, 14).Value = 0
$ws.Cells.Item(127, 15).Value = 2
$ws.Cells.Item(127, 16).Value = 2
Set-TextCell 127 17 "41%"
Set-TextCell 127 18 "59%"
Set-TextCell 127 19 "E"

# Row 128: Huachipato vs O'Higgins
Set-TextCell 128 1 "2025-07-21"
Set-TextCell 128 2 "Huachipato"
Set-TextCell 128 3 "O'Higgins"
$ws.Cells.Item(128, 4).Value = 2
$ws.Cells.Item(128, 5).Value = 1
$ws.Cells.Item(128, 6).Value = 1339133
$ws.Cells.Item(128, 7).Value = 6
$ws.Cells.Item(128, 8).Value = 4
$ws.Cells.Item(128, 9).Value = 2
$ws.Cells.Item(128, 10).Value = 2
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 1
$ws.Cells.Item(128, 13).Value = 0
$ws.Cells.Item(128, 14).Value = 0
$ws.Cells.Item(128, 15).Value = 2
$ws.Cells.Item(128, 16).Value = 1
Set-TextCell 128 17 "52%"
Set-TextCell 128 18 "48%"
Set-TextCell 128 19 "L"

Write-Output "Added rows 122-128"